$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickets = @{
    2 = 378918
    3 = 835418
    4 = 832965
    5 = 385219
    6 = 178002
    7 = 336085
    8 = 336491
    9 = 645991
}

foreach ($row in $tickets.Keys) {
    $id = $tickets[$row]
    $ws.Cells.Item($row, 4).Value = $id
    $ws.Cells.Item($row, 5).Value = "https://93mtzf.deta.dev/ticket/$id"
}

$ws.Range("A10:E11").Delete()
